$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 82
$ws.Cells.Item($row, 1).Value = "2024-10-23 00:00:00"
$ws.Cells.Item($row, 2).Value = 73650
$ws.Cells.Item($row, 3).Value = 10305.74
$ws.Cells.Item($row, 4).Value = 9120.129999999999
$ws.Cells.Item($row, 5).Value = 7.1285
